$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Running only search suite: set Runmode to "N" for every suite except
# the Search module (B Suite, row 3), which keeps "Y".
$ws.Range("C2").Value = "N"
$ws.Range("C4").Value = "N"
$ws.Range("C5").Value = "N"
$ws.Range("C6").Value = "N"
$ws.Range("C7").Value = "N"

# Update the sheet selection to reflect the last edited cell.
$ws.Range("C7").Select()
